$d = $word.ActiveDocument

# --- 1) Fix the typo "Capaidade: " -> "Capacidade: " -------------------
# Locate the point right after "Capa" (before "idade: ") and type the
# missing "c", exactly like a user clicking there and typing one letter.
$rng = $d.Content
$rng.Find.Execute("Capa", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$rng.InsertAfter("c")

# Word drops a "_GoBack" bookmark at the position of the last edit.
$afterC = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $afterC)

# --- 2) Tag the picture run's language (eastAsia = pt-BR) --------------
$p2 = $d.Paragraphs(2)
$p2.Range.LanguageIDFarEast = "pt-BR"
